$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.830.44"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "2.496.79"
$ws.Range("E3").Value = "  -3.88%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.87"
$ws.Range("E5").Value = "  -4.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.59"
$ws.Range("E6").Value = "  -5.19%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("D9").Value = "2.496.44"
$ws.Range("E9").Value = "  -3.81%  "
$ws.Range("E10").Value = "  -9.53%  "
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("E12").Value = "  -7.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  -6.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.07"
$ws.Range("E14").Value = "  -6.90%  "
$ws.Range("D15").Value = "2.947.63"
$ws.Range("E15").Value = "  -3.84%  "
$ws.Range("D16").Value = "61.744.94"
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000162"
$ws.Range("E17").Value = "  -8.33%  "
$ws.Range("D18").Value = "2.496.22"
$ws.Range("E18").Value = "  -3.57%  "
$ws.Range("E19").Value = "  -7.20%  "
$ws.Range("E20").Value = "  -6.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.18"
$ws.Range("E21").Value = "  -7.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "320.77"
$ws.Range("E22").Value = "  -5.92%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.66"
$ws.Range("E24").Value = "  -5.61%  "
$ws.Range("E25").Value = "  -3.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000100"
$ws.Range("E26").Value = "  -6.88%  "
$ws.Range("D27").Value = "2.622.14"
$ws.Range("E27").Value = "  -3.30%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  -5.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.37"
$ws.Range("E30").Value = "  -8.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "533.84"
$ws.Range("E31").Value = "  -7.48%  "
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("E33").Value = "  -6.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("E34").Value = "  -8.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  -9.41%  "
$ws.Range("E36").Value = "  -9.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.87"
$ws.Range("E37").Value = "  -9.49%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("E39").Value = "  -5.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.47"
$ws.Range("E40").Value = "  -6.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "144.03"
$ws.Range("E41").Value = "  -6.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  -9.23%  "
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.29"
$ws.Range("E45").Value = "  -7.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "148.69"
$ws.Range("E46").Value = "  -5.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.57"
$ws.Range("E47").Value = "  -8.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.75"
$ws.Range("E48").Value = "  -10.68%  "
$ws.Range("E49").Value = "  -8.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.588"
$ws.Range("E50").Value = "  -6.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0942"
$ws.Range("E51").Value = "  -5.76%  "
